# Weekly refresh: insert a new daily price record for
# "Femacal de La Calera - Pepino ensalada" at row 650, pushing the
# existing rows 650-683 down to 651-684 (their data is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at 650 (shifts 650..683 down to 651..684).
$ws.Rows(650).Insert()

# Populate the new row 650 with the latest price observation.
$ws.Cells.Item(650, 1).Value  = 3
$ws.Cells.Item(650, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(650, 3).Value  = "Coquimbo"
$ws.Cells.Item(650, 4).Value  = 45267
$ws.Cells.Item(650, 5).Value  = 5
$ws.Cells.Item(650, 6).Value  = 100112043
$ws.Cells.Item(650, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(650, 8).Value  = "Sin especificar"
$ws.Cells.Item(650, 9).Value  = "Primera"
$ws.Cells.Item(650, 10).Value = 95
$ws.Cells.Item(650, 11).Value = 18000
$ws.Cells.Item(650, 12).Value = 19000
$ws.Cells.Item(650, 13).Value = 18526
$ws.Cells.Item(650, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(650, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(650, 16).Value = 265
$ws.Cells.Item(650, 17).Value = 70
$ws.Cells.Item(650, 18).Value = "Hortaliza"
